$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.811.00"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.448.73"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").Value = "2.445.42"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  -3.98%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").Value = "2.871.59"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "62.617.76"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "2.440.12"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  -3.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "640.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.82%  "
$ws.Range("D28").Value = "0.0₃0973"
$ws.Range("E28").Value = "  -6.26%  "
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "148.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.601"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "0.0₆0238"
$ws.Range("E51").Value = "  +8.67%  "
